$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Test data was re-shuffled: the old row 6 (debashree.p@insync.co.inxxxxxxx /
# Efgh.1234) moved down to row 7 and row 6 itself became a blank spacer row;
# two brand-new rows were appended (row 8: rio1@yopmail.com, row 9: blank /
# abcd.1234).
# ---------------------------------------------------------------------------

# 1) Write the new cell values first (hyperlink text + plain text).
$ws.Range("A7").Value = "debashree.p@insync.co.inxxxxxxx"
$ws.Range("B7").Value = "Efgh.1234"

$ws.Range("A8").Value = "rio1@yopmail.com"

$ws.Range("B9").Value = "abcd.1234"

# 2) Rebuild the hyperlinks. (The host's hyperlink collection only supports
#    deleting the whole sheet collection at once, so delete everything and
#    re-add all six links - the four that are unchanged plus the two that
#    moved/are new - in document order so the relationship ids line up.)
$ws.Range("A6").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:debashree.p@insync.co.in") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:debashree.p@insync.co.in") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:rio1@yopmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:debashree.p@insync.co.in") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:debashree.p@insync.co.inxxxxxxx") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:rio1@yopmail.com") | Out-Null

# 3) Adding hyperlinks re-stamps a cell's style, so restore the original
#    "Hyperlink" formatting (bordered + underlined) on every A-column link
#    cell, and the plain bordered formatting on every B-column cell / the
#    non-link A9 cell, by pasting formats from already-correctly-styled
#    cells (this reuses the existing style indexes instead of creating new
#    duplicate ones).
$ws.Range("A2").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B9").PasteSpecial(-4122)

# 4) Row 6 becomes an empty spacer row (keeps its formatting, loses content).
$ws.Range("A6").ClearContents()
$ws.Range("B6").ClearContents()

# 5) Match the saved selection/view state.
$ws.Range("A1:B9").Select() | Out-Null
